# 10busCDF.xlsx edit: add a 10th bus to the dataset and convert the
# load columns (H/I) from per-unit fractions to raw MW/MVAR-like
# magnitudes for testing purposes (see commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Update the existing BUS DATA load columns (H = col 8, I = col 9)
#    for buses 2-9 (rows 4-11) to the new non-PU magnitudes.
# ---------------------------------------------------------------
$busLoads = @{
    4  = @(1500, 750)
    5  = @(900, 500)
    6  = @(800, 450)
    7  = @(700, 400)
    8  = @(600, 350)
    9  = @(500, 300)
    10 = @(400, 250)
    11 = @(300, 200)
}

foreach ($row in $busLoads.Keys) {
    $vals = $busLoads[$row]
    $ws.Cells.Item($row, 8).Value = $vals[0]
    $ws.Cells.Item($row, 9).Value = $vals[1]
}

# ---------------------------------------------------------------
# 2) Insert a new row at row 12 (right after the last existing bus,
#    bus 9) to hold the data for the new "Bus 10". Inserting a whole
#    row here shifts everything below (the -999 marker, branch data,
#    loss zones, tie lines, etc.) down by one row, which matches the
#    target layout.
# ---------------------------------------------------------------
$ws.Rows(12).Insert()

$ws.Cells.Item(12, 1).Value  = 10
$ws.Cells.Item(12, 2).Value  = "Bus 10"
$ws.Cells.Item(12, 3).Value  = "MV 1"
$ws.Cells.Item(12, 4).Value  = 1
$ws.Cells.Item(12, 5).Value  = 0
$ws.Cells.Item(12, 6).Value  = 0.9345
$ws.Cells.Item(12, 7).Value  = -0.1222
$ws.Cells.Item(12, 8).Value  = 200
$ws.Cells.Item(12, 9).Value  = 100
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 12.66
$ws.Cells.Item(12, 13).Value = 0
$ws.Cells.Item(12, 14).Value = 0
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 0
$ws.Cells.Item(12, 18).Value = 0

# ---------------------------------------------------------------
# 3) Update the BRANCH DATA impedance columns (G = col 7, H = col 8)
#    for the nine branches that now live in rows 15-23 (they were
#    rows 14-22 before the row-12 insertion shifted them down).
# ---------------------------------------------------------------
$branchImpedances = @{
    15 = @(0.1705, 0.3409)
    16 = @(0.2273, 0.4545)
    17 = @(0.2273, 0.4545)
    18 = @(0.2273, 0.4545)
    19 = @(0.2273, 0.4545)
    20 = @(0.2273, 0.4545)
    21 = @(0.2273, 0.4545)
    22 = @(0.2273, 0.4545)
    23 = @(0.2273, 0.4545)
}

foreach ($row in $branchImpedances.Keys) {
    $vals = $branchImpedances[$row]
    $ws.Cells.Item($row, 7).Value = $vals[0]
    $ws.Cells.Item($row, 8).Value = $vals[1]
}

# ---------------------------------------------------------------
# 4) Restore the active cell selection to H8 (matching the saved
#    view state in the edited workbook).
# ---------------------------------------------------------------
[void]$ws.Range("H8").Select()
